# Derek's Log - add Friday (06/16/2016-serial 42659) entries to the Logs sheet
# and correct two AV Shutdown times on the preceding day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$ws.Activate()

# ---------------------------------------------------------------
# 1) Fix the two "2200" -> "2150" AV Shutdown times on 2016-12-06
# ---------------------------------------------------------------
$ws.Cells.Item(352, 3).Value = "2150"
$ws.Cells.Item(355, 3).Value = "2150"

# ---------------------------------------------------------------
# 2) Append the new Friday block: a day-separator row (356) plus
#    six data rows (357-362). Formats are copied, cell by cell,
#    from existing rows that already carry the exact style we
#    need, then the values are overwritten for the new entries.
# ---------------------------------------------------------------

# Row 356: day-separator row ("FRIDAY"), same look as row 41.
$ws.Range("A41:F41").Copy()
$ws.Range("A356:F356").PasteSpecial(-4122)
$ws.Cells.Item(356, 3).Value = "FRIDAY"

# Row 357: AV Shutdown 1530 / R / S203 (no notes) - same look as row 197.
$ws.Range("A197:E197").Copy()
$ws.Range("A357:E357").PasteSpecial(-4122)
$ws.Cells.Item(357, 1).Value = "AV Shutdown"
$ws.Cells.Item(357, 2).Value = 42659
$ws.Cells.Item(357, 3).Value = "1530"
$ws.Cells.Item(357, 4).Value = "R"
$ws.Cells.Item(357, 5).Value = "S203"

# Row 358: AV Shutdown 1530 / R / N203 (no notes) - same look as row 324.
$ws.Range("A324:E324").Copy()
$ws.Range("A358:E358").PasteSpecial(-4122)
$ws.Cells.Item(358, 1).Value = "AV Shutdown"
$ws.Cells.Item(358, 2).Value = 42659
$ws.Cells.Item(358, 3).Value = "1530"
$ws.Cells.Item(358, 4).Value = "R"
$ws.Cells.Item(358, 5).Value = "N203"

# Row 359: Pickup PC 1600 / VC / 256 / DLP TV note - same look as row 346.
$ws.Range("A346:F346").Copy()
$ws.Range("A359:F359").PasteSpecial(-4122)
$ws.Cells.Item(359, 1).Value = "Pickup PC"
$ws.Cells.Item(359, 2).Value = 42659
$ws.Cells.Item(359, 3).Value = "1600"
$ws.Cells.Item(359, 4).Value = "VC"
$ws.Cells.Item(359, 5).Value = "256"
$ws.Cells.Item(359, 6).Value = "Return large screen DLP TV to Vanier 132 storeroom. Pick up wireless keyboard and remote control. "

# Row 360: AV Shutdown 1800 / R / N102 / Nat Taylor Cinema note - same look as row 352.
$ws.Range("A352:F352").Copy()
$ws.Range("A360:F360").PasteSpecial(-4122)
$ws.Cells.Item(360, 1).Value = "AV Shutdown"
$ws.Cells.Item(360, 2).Value = 42659
$ws.Cells.Item(360, 3).Value = "1800"
$ws.Cells.Item(360, 4).Value = "R"
$ws.Cells.Item(360, 5).Value = "N102"
$ws.Cells.Item(360, 6).Value = "Nat Taylor Cinema. Lock cinema all doors after shutdown."

# Row 361: Other 1730 / MC / 157A / Door code note - same look as row 325.
$ws.Range("A325:F325").Copy()
$ws.Range("A361:F361").PasteSpecial(-4122)
$ws.Cells.Item(361, 1).Value = "Other"
$ws.Cells.Item(361, 2).Value = 42659
$ws.Cells.Item(361, 3).Value = "1730"
$ws.Cells.Item(361, 4).Value = "MC"
$ws.Cells.Item(361, 5).Value = "157A"
$ws.Cells.Item(361, 6).Value = "Door code 11012* "

# Row 362: AV Shutdown 2050 / FC / 152 / Founders Assembly Hall note - same look as row 48.
$ws.Range("A48:F48").Copy()
$ws.Range("A362:F362").PasteSpecial(-4122)
$ws.Cells.Item(362, 1).Value = "AV Shutdown"
$ws.Cells.Item(362, 2).Value = 42659
$ws.Cells.Item(362, 3).Value = "2050"
$ws.Cells.Item(362, 4).Value = "FC"
$ws.Cells.Item(362, 5).Value = "152"
$ws.Cells.Item(362, 6).Value = "Founders Assembly Hall - group in room using projector and computer in room. No order - please make sure projector gets turned off (remote in Fdrs 156A storeroom. Log off PC and please LOCK ROOM. Key for room in Fdrs 164 storeroom."

# Row heights that differ from the sheet default (matches author's manual sizing).
$ws.Rows.Item(359).RowHeight = 30
$ws.Rows.Item(362).RowHeight = 60

# ---------------------------------------------------------------
# 3) Restore the selection to match where the author ended up.
# ---------------------------------------------------------------
$ws.Range("B358:B362").Select()
